$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, whether it is a Price/Volume column that
# must be forced to Text so Excel does not auto-convert numeric-looking strings
# (e.g. "503.18") into real numbers.
$updates = @(
    @{Cell='D2'; Value='54.527.08'; Text=$true}
    @{Cell='E2'; Value='  +0.39%  '; Text=$true}
    @{Cell='D3'; Value='2.275.08'; Text=$true}
    @{Cell='E3'; Value='  -0.15%  '; Text=$true}
    @{Cell='E4'; Value='  +0.04%  '; Text=$true}
    @{Cell='D5'; Value='503.18'; Text=$true}
    @{Cell='E5'; Value='  +0.92%  '; Text=$true}
    @{Cell='D6'; Value='128.38'; Text=$true}
    @{Cell='E6'; Value='  +0.07%  '; Text=$true}
    @{Cell='D7'; Value='0.996'; Text=$true}
    @{Cell='E7'; Value='  -0.19%  '; Text=$true}
    @{Cell='D8'; Value='0.528'; Text=$true}
    @{Cell='E8'; Value='  -0.25%  '; Text=$true}
    @{Cell='D9'; Value='2.287.73'; Text=$true}
    @{Cell='E9'; Value='  +0.10%  '; Text=$true}
    @{Cell='D10'; Value='0.0963'; Text=$true}
    @{Cell='E10'; Value='  +0.83%  '; Text=$true}
    @{Cell='E11'; Value='  +1.30%  '; Text=$true}
    @{Cell='E12'; Value='  +1.88%  '; Text=$true}
    @{Cell='D13'; Value='4.89'; Text=$true}
    @{Cell='E13'; Value='  +3.79%  '; Text=$true}
    @{Cell='D14'; Value='23.31'; Text=$true}
    @{Cell='E14'; Value='  +3.16%  '; Text=$true}
    @{Cell='D15'; Value='2.679.05'; Text=$true}
    @{Cell='E15'; Value='  -0.15%  '; Text=$true}
    @{Cell='D16'; Value='54.590.52'; Text=$true}
    @{Cell='E16'; Value='  +0.63%  '; Text=$true}
    @{Cell='E17'; Value='  +0.79%  '; Text=$true}
    @{Cell='D18'; Value='2.283.69'; Text=$true}
    @{Cell='E18'; Value='  +0.15%  '; Text=$true}
    @{Cell='D19'; Value='10.33'; Text=$true}
    @{Cell='E19'; Value='  +0.44%  '; Text=$true}
    @{Cell='D20'; Value='4.13'; Text=$true}
    @{Cell='E20'; Value='  +0.20%  '; Text=$true}
    @{Cell='D21'; Value='306.54'; Text=$true}
    @{Cell='E21'; Value='  +0.66%  '; Text=$true}
    @{Cell='D22'; Value='6.43'; Text=$true}
    @{Cell='E22'; Value='  +0.12%  '; Text=$true}
    @{Cell='D23'; Value='1.00'; Text=$true}
    @{Cell='E23'; Value='  +0.14%  '; Text=$true}
    @{Cell='D24'; Value='60.29'; Text=$true}
    @{Cell='E24'; Value='  -2.51%  '; Text=$true}
    @{Cell='D25'; Value='0.997'; Text=$true}
    @{Cell='E25'; Value='  -0.17%  '; Text=$true}
    @{Cell='E26'; Value='  -1.19%  '; Text=$true}
    @{Cell='D27'; Value='7.44'; Text=$true}
    @{Cell='E27'; Value='  +1.56%  '; Text=$true}
    @{Cell='D28'; Value='170.85'; Text=$true}
    @{Cell='E28'; Value='  -2.92%  '; Text=$true}
    @{Cell='B29'; Value='PEPE'; Text=$false}
    @{Cell='C29'; Value='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; Text=$false}
    @{Cell='D29'; Value='0.0₃0701'; Text=$true}
    @{Cell='E29'; Value='  +2.20%  '; Text=$true}
    @{Cell='B30'; Value='Aptos'; Text=$false}
    @{Cell='C30'; Value='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; Text=$false}
    @{Cell='D30'; Value='6.05'; Text=$true}
    @{Cell='E30'; Value='  +1.58%  '; Text=$true}
    @{Cell='B31'; Value='PancakeSwap'; Text=$false}
    @{Cell='C31'; Value='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; Text=$false}
    @{Cell='D31'; Value='1.62'; Text=$true}
    @{Cell='E31'; Value='  +0.41%  '; Text=$true}
    @{Cell='E32'; Value='  +3.06%  '; Text=$true}
    @{Cell='E33'; Value='  +0.01%  '; Text=$true}
    @{Cell='D34'; Value='17.91'; Text=$true}
    @{Cell='E34'; Value='  +0.83%  '; Text=$true}
    @{Cell='D35'; Value='0.995'; Text=$true}
    @{Cell='D36'; Value='0.908'; Text=$true}
    @{Cell='E36'; Value='  -1.83%  '; Text=$true}
    @{Cell='E37'; Value='  +0.11%  '; Text=$true}
    @{Cell='D38'; Value='3.77'; Text=$true}
    @{Cell='E38'; Value='  +0.67%  '; Text=$true}
    @{Cell='D39'; Value='36.46'; Text=$true}
    @{Cell='E39'; Value='  +0.73%  '; Text=$true}
    @{Cell='D40'; Value='0.374'; Text=$true}
    @{Cell='E40'; Value='  +0.08%  '; Text=$true}
    @{Cell='E41'; Value='  +0.34%  '; Text=$true}
    @{Cell='D42'; Value='5.04'; Text=$true}
    @{Cell='E42'; Value='  +5.70%  '; Text=$true}
    @{Cell='E43'; Value='  +0.04%  '; Text=$true}
    @{Cell='D44'; Value='126.56'; Text=$true}
    @{Cell='E44'; Value='  +0.63%  '; Text=$true}
    @{Cell='D45'; Value='0.0496'; Text=$true}
    @{Cell='E45'; Value='  +1.09%  '; Text=$true}
    @{Cell='D46'; Value='248.95'; Text=$true}
    @{Cell='E46'; Value='  +3.68%  '; Text=$true}
    @{Cell='D47'; Value='0.0903'; Text=$true}
    @{Cell='E47'; Value='  +0.67%  '; Text=$true}
    @{Cell='E48'; Value='  +0.02%  '; Text=$true}
    @{Cell='E49'; Value='  +0.27%  '; Text=$true}
    @{Cell='D50'; Value='0.0207'; Text=$true}
    @{Cell='E50'; Value='  +0.47%  '; Text=$true}
    @{Cell='D51'; Value='10.81'; Text=$true}
    @{Cell='E51'; Value='  +0.36%  '; Text=$true}
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    if ($u.Text) {
        $c.NumberFormat = "@"
        $c.Value = $u.Value
        $c.Style = "Normal"
    } else {
        $c.Value = $u.Value
    }
}
